$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Refresh the scraped crypto price/volume figures (and the re-ranked rows
# 13-15 / 38-39, whose Coin/Link/Price/Volume shifted down a slot).
# Values are prefixed with a leading apostrophe so Excel keeps them as
# literal text (matching the original inlineStr cells) instead of
# auto-coercing numeric-looking strings like "0.9988" into numbers.
$ws.Range('D2').Value = "'29.408.35"
$ws.Range('E2').Value = "'  -0.36%  "
$ws.Range('D3').Value = "'1.847.55"
$ws.Range('E3').Value = "'  -0.16%  "
$ws.Range('D4').Value = "'0.9988"
$ws.Range('E4').Value = "'  +0.05%  "
$ws.Range('D5').Value = "'240.77"
$ws.Range('E5').Value = "'  -0.97%  "
$ws.Range('D6').Value = "'0.6317"
$ws.Range('E6').Value = "'  -3.35%  "
$ws.Range('D7').Value = "'0.9998"
$ws.Range('E7').Value = "'  +0.05%  "
$ws.Range('D8').Value = "'0.07594"
$ws.Range('E8').Value = "'  +1.16%  "
$ws.Range('E9').Value = "'  -0.30%  "
$ws.Range('D10').Value = "'24.49"
$ws.Range('E10').Value = "'  -0.13%  "
$ws.Range('D11').Value = "'2.386.85"
$ws.Range('E11').Value = "'  +28.60%  "
$ws.Range('D12').Value = "'0.07723"
$ws.Range('E12').Value = "'  +1.22%  "
$ws.Range('B13').Value = "'WrappedliquidstakedEther2.0"
$ws.Range('C13').Value = "'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
$ws.Range('D13').Value = "'2.514.76"
$ws.Range('E13').Value = "'  +18.23%  "
$ws.Range('B14').Value = "'Polkadot"
$ws.Range('C14').Value = "'https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range('D14').Value = "'4.986"
$ws.Range('E14').Value = "'  -0.75%  "
$ws.Range('B15').Value = "'Polygon"
$ws.Range('C15').Value = "'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
$ws.Range('D15').Value = "'0.6863"
$ws.Range('E15').Value = "'  +0.13%  "
$ws.Range('D16').Value = "'83.00"
$ws.Range('E16').Value = "'  -0.87%  "
$ws.Range('D17').Value = "'0.000009914"
$ws.Range('E17').Value = "'  +4.32%  "
$ws.Range('D18').Value = "'6.169"
$ws.Range('E18').Value = "'  +0.80%  "
$ws.Range('D19').Value = "'29.444.52"
$ws.Range('E19').Value = "'  -0.27%  "
$ws.Range('D20').Value = "'231.78"
$ws.Range('E20').Value = "'  -2.46%  "
$ws.Range('E21').Value = "'  -0.67%  "
$ws.Range('D22').Value = "'0.9996"
$ws.Range('E22').Value = "'  +0.01%  "
$ws.Range('D23').Value = "'7.610"
$ws.Range('E23').Value = "'  -1.28%  "
$ws.Range('E24').Value = "'  +0.01%  "
$ws.Range('D25').Value = "'154.28"
$ws.Range('E25').Value = "'  -1.61%  "
$ws.Range('E26').Value = "'  -2.04%  "
$ws.Range('E27').Value = "'  -0.35%  "
$ws.Range('E28').Value = "'  -0.73%  "
$ws.Range('D29').Value = "'1.471"
$ws.Range('E29').Value = "'  -1.01%  "
$ws.Range('E30').Value = "'  -3.97%  "
$ws.Range('D31').Value = "'1.257"
$ws.Range('E31').Value = "'  +1.06%  "
$ws.Range('E32').Value = "'  -0.21%  "
$ws.Range('D33').Value = "'4.028"
$ws.Range('E33').Value = "'  -1.30%  "
$ws.Range('D34').Value = "'1.865"
$ws.Range('E34').Value = "'  +0.45%  "
$ws.Range('E35').Value = "'  -1.99%  "
$ws.Range('E36').Value = "'  -0.85%  "
$ws.Range('D37').Value = "'2.592"
$ws.Range('E37').Value = "'  -0.10%  "
$ws.Range('B38').Value = "'RocketPoolETH"
$ws.Range('C38').Value = "'https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth"
$ws.Range('D38').Value = "'2.447.88"
$ws.Range('E38').Value = "'  +20.38%  "
$ws.Range('B39').Value = "'Maker"
$ws.Range('C39').Value = "'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range('D39').Value = "'1.246.77"
$ws.Range('E39').Value = "'  +3.74%  "
$ws.Range('D40').Value = "'2.793"
$ws.Range('E40').Value = "'  -0.27%  "
$ws.Range('D41').Value = "'0.01807"
$ws.Range('E41').Value = "'  +1.30%  "
$ws.Range('D42').Value = "'0.9067"
$ws.Range('E42').Value = "'  -0.09%  "
$ws.Range('D43').Value = "'6.116"
$ws.Range('E43').Value = "'  -2.12%  "
$ws.Range('E44').Value = "'  +0.00%  "
$ws.Range('D45').Value = "'67.27"
$ws.Range('E45').Value = "'  +1.05%  "
$ws.Range('D46').Value = "'101.41"
$ws.Range('E46').Value = "'  -0.46%  "
$ws.Range('D47').Value = "'7.309"
$ws.Range('E47').Value = "'  -2.04%  "
$ws.Range('E48').Value = "'  +0.59%  "
$ws.Range('D49').Value = "'9.177"
$ws.Range('E49').Value = "'  +0.45%  "
$ws.Range('D50').Value = "'0.4014"
$ws.Range('E50').Value = "'  -1.04%  "
$ws.Range('D51').Value = "'1.697"
$ws.Range('E51').Value = "'  +2.13%  "
